$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the placeholder "Capital Call ..." Call Name values with the real
# call names for this Fund Unit upload.
$ws.Range("C2:C4").Value = "Call 1"
$ws.Range("C5:C6").Value = "Call 2"

# Remove the stray tiny Arial font formatting that was on column C, restoring
# the default cell style.
$ws.Range("C2:C6").ClearFormats() | Out-Null

# Leave the selection where the user finished editing.
$ws.Range("C7").Select() | Out-Null
